# Update SPR valid example: append rows 3-8 to the "Dataset" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

# Trailing "'" (a bare quote-prefix) writes an empty, but present, TEXT
# cell -- matching the workbook's existing empty "Comment" cells (K column),
# which are stored as empty inline/shared strings rather than blank cells.
$data = @(
    @("COVIC 1",  "Nucleoprotein 1", "4", ">500000",  "NA",   "NA",   "NA",   "23100000",   "14500", "negative", "'"),
    @("COVIC 1",  "Spike protein 2", "6", "32100",    "6100", "14.6", "5",    "10",         "0.01",  "positive", "'"),
    @("COVIC 4",  "Nucleoprotein 1", "4", ">500000",  "NA",   "0.1",  "0.03", ">200000000", "NA",    "negative", "'"),
    @("COVIC 4",  "Spike protein 1", "6", "674000",   "10000","347",  "12",   "12.3",       "0.01",  "positive", "'"),
    @("COVIC 10", "Spike protein 2", "6", ">500000",  "NA",   "0.1",  "0.03", ">200000000", "NA",    "negative", "'"),
    @("COVIC 10", "Nucleoprotein 1", "4", "82000",    "4000", "29",   "4.1",  "<.1",        "0.01",  "positive", "'")
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $rng = $ws.Range("A" + $r + ":K" + $r)
    # Force text storage so numeric-looking strings (e.g. "6", "491000")
    # stay text, matching the rest of the sheet.
    $rng.NumberFormat = "@"
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $cell.Value = $rowVals[$c]
    }
    # Drop the temporary "Text" number-format style so the new cells end up
    # unstyled, like the existing data rows.
    $rng.ClearFormats()
}
